$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E21").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F21").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G21").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H21").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I21").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J21").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K21").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L21").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M21").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N21").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E35").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F35").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G35").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H35").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I35").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J35").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K35").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L35").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M35").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N35").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E49").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F49").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G49").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H49").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I49").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J49").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K49").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L49").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M49").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N49").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E62").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F62").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G62").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H62").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I62").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J62").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K62").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L62").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M62").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N62").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E75").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F75").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G75").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H75").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I75").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J75").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K75").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L75").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M75").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N75").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 241
$ws.Range("H12").Value = 0
$ws.Range("E14").Value = 2323
$ws.Range("F14").Value = 2407
$ws.Range("G14").Value = 2134
$ws.Range("H14").Value = 2062
$ws.Range("I14").Value = 2113
$ws.Range("J14").Value = 2030
$ws.Range("K14").Value = 2020
$ws.Range("L14").Value = 1447
$ws.Range("M14").Value = 2101
$ws.Range("N14").Value = 1968
$ws.Range("E15").Value = 4543
$ws.Range("F15").Value = 4504
$ws.Range("G15").Value = 4275
$ws.Range("H15").Value = 4106
$ws.Range("I15").Value = 4755
$ws.Range("J15").Value = 4707
$ws.Range("K15").Value = 4695
$ws.Range("L15").Value = 4989
$ws.Range("M15").Value = 5014
$ws.Range("N15").Value = 4749
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 55
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 71
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 109
$ws.Range("K16").Value = 30
$ws.Range("L16").Value = 72
$ws.Range("M16").Value = 10
$ws.Range("N16").Value = 129
$ws.Range("E17").Value = 6867
$ws.Range("F17").Value = 6966
$ws.Range("G17").Value = 6443
$ws.Range("H17").Value = 6239
$ws.Range("I17").Value = 6871
$ws.Range("J17").Value = 6846
$ws.Range("K17").Value = 6745
$ws.Range("L17").Value = 6508
$ws.Range("M17").Value = 7125
$ws.Range("N17").Value = 7087
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 1
$ws.Range("E24").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = -6
$ws.Range("G25").Value = "-"
$ws.Range("J26").Value = "-"
$ws.Range("K26").Value = 0
$ws.Range("E28").Value = 2158
$ws.Range("F28").Value = 2676
$ws.Range("G28").Value = 1916
$ws.Range("H28").Value = 2063
$ws.Range("I28").Value = 2064
$ws.Range("J28").Value = "-"
$ws.Range("K28").Value = 1793
$ws.Range("L28").Value = 2031
$ws.Range("M28").Value = 2233
$ws.Range("N28").Value = 2001
$ws.Range("E29").Value = 4514
$ws.Range("F29").Value = 4520
$ws.Range("G29").Value = 4145
$ws.Range("H29").Value = 4123
$ws.Range("I29").Value = 4686
$ws.Range("J29").Value = "-"
$ws.Range("K29").Value = 4755
$ws.Range("L29").Value = 4841
$ws.Range("M29").Value = 5083
$ws.Range("N29").Value = 4845
$ws.Range("E30").Value = 42
$ws.Range("F30").Value = 45
$ws.Range("H30").Value = 32
$ws.Range("I30").Value = 41
$ws.Range("J30").Value = "-"
$ws.Range("K30").Value = 33
$ws.Range("L30").Value = 25
$ws.Range("M30").Value = 59
$ws.Range("N30").Value = 64
$ws.Range("E31").Value = 6714
$ws.Range("F31").Value = 7235
$ws.Range("G31").Value = 6106
$ws.Range("H31").Value = 6218
$ws.Range("I31").Value = 6791
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6581
$ws.Range("L31").Value = 6897
$ws.Range("M31").Value = 7375
$ws.Range("N31").Value = 6911
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 1946
$ws.Range("E38").Value = "-"
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = -4151
$ws.Range("G39").Value = "-"
$ws.Range("E40").Value = -79043
$ws.Range("F40").Value = -120553
$ws.Range("G40").Value = -28419
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = "-"
$ws.Range("K40").Value = 0
$ws.Range("E42").Value = 996588
$ws.Range("F42").Value = 1616599
$ws.Range("G42").Value = 1083544
$ws.Range("H42").Value = 1149756
$ws.Range("I42").Value = 1261918
$ws.Range("J42").Value = "-"
$ws.Range("K42").Value = 1253660
$ws.Range("L42").Value = 1543131
$ws.Range("M42").Value = 1847502
$ws.Range("N42").Value = 1786300
$ws.Range("E43").Value = 2047327
$ws.Range("F43").Value = 2598576
$ws.Range("G43").Value = 2288232
$ws.Range("H43").Value = 2248458
$ws.Range("I43").Value = 2559915
$ws.Range("J43").Value = "-"
$ws.Range("K43").Value = 2766238
$ws.Range("L43").Value = 3490127
$ws.Range("M43").Value = 3818889
$ws.Range("N43").Value = 3848163
$ws.Range("E44").Value = 15368
$ws.Range("F44").Value = 18278
$ws.Range("G44").Value = 15272
$ws.Range("H44").Value = 10421
$ws.Range("I44").Value = 13656
$ws.Range("J44").Value = "-"
$ws.Range("K44").Value = 16598
$ws.Range("L44").Value = 16891
$ws.Range("M44").Value = 43391
$ws.Range("N44").Value = 45729
$ws.Range("E45").Value = 2980240
$ws.Range("F45").Value = 4108749
$ws.Range("G45").Value = 3358629
$ws.Range("H45").Value = 3408635
$ws.Range("I45").Value = 3835489
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4036496
$ws.Range("L45").Value = 5050149
$ws.Range("M45").Value = 5709782
$ws.Range("N45").Value = 5682138
$ws.Range("N51").Value = 1946000000
$ws.Range("E53").Value = "-"
$ws.Range("F53").Value = 691833333
$ws.Range("G53").Value = "-"
$ws.Range("E56").Value = 446714219
$ws.Range("F56").Value = 581588015
$ws.Range("G56").Value = 565524008
$ws.Range("H56").Value = 545851672
$ws.Range("I56").Value = 608153253
$ws.Range("J56").Value = 670054780
$ws.Range("K56").Value = 699196877
$ws.Range("L56").Value = 759788774
$ws.Range("M56").Value = 827363189
$ws.Range("N56").Value = 892703648
$ws.Range("E57").Value = 451400975
$ws.Range("F57").Value = 571360619
$ws.Range("G57").Value = 552046321
$ws.Range("H57").Value = 544193791
$ws.Range("I57").Value = 546290013
$ws.Range("J57").Value = 550776341
$ws.Range("K57").Value = 581753523
$ws.Range("L57").Value = 720951663
$ws.Range("M57").Value = 751306118
$ws.Range("N57").Value = 794254489
$ws.Range("E58").Value = 365690476
$ws.Range("F58").Value = 406133333
$ws.Range("G58").Value = 339377778
$ws.Range("H58").Value = 325406250
$ws.Range("I58").Value = 333073171
$ws.Range("J58").Value = 338758621
$ws.Range("K58").Value = 502969697
$ws.Range("L58").Value = 675640000
$ws.Range("M58").Value = 735440678
$ws.Range("N58").Value = 714515625
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = -619
$ws.Range("E65").Value = "-"
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = 1622
$ws.Range("G66").Value = "-"
$ws.Range("E68").Value = -672841
$ws.Range("F68").Value = -825922
$ws.Range("G68").Value = -889611
$ws.Range("H68").Value = -816013
$ws.Range("I68").Value = -1028232
$ws.Range("J68").Value = -1023295
$ws.Range("K68").Value = -1006339
$ws.Range("L68").Value = -1306718
$ws.Range("M68").Value = -1433154
$ws.Range("N68").Value = -1358056
$ws.Range("E69").Value = -1269475
$ws.Range("F69").Value = -1480421
$ws.Range("G69").Value = -1619292
$ws.Range("H69").Value = -2131962
$ws.Range("I69").Value = -2311624
$ws.Range("J69").Value = -2660595
$ws.Range("K69").Value = -2823893
$ws.Range("L69").Value = -3248247
$ws.Range("M69").Value = -3404487
$ws.Range("N69").Value = -3251351
$ws.Range("E70").Value = -10459
$ws.Range("F70").Value = -10045
$ws.Range("G70").Value = -13227
$ws.Range("H70").Value = -8634
$ws.Range("I70").Value = -13574
$ws.Range("J70").Value = -11843
$ws.Range("K70").Value = -13541
$ws.Range("L70").Value = -13276
$ws.Range("M70").Value = -32680
$ws.Range("N70").Value = -6191
$ws.Range("E71").Value = -1952775
$ws.Range("F71").Value = -2314766
$ws.Range("G71").Value = -2522130
$ws.Range("H71").Value = -2956609
$ws.Range("I71").Value = -3353430
$ws.Range("J71").Value = -3695733
$ws.Range("K71").Value = -3843773
$ws.Range("L71").Value = -4568241
$ws.Range("M71").Value = -4870321
$ws.Range("N71").Value = -4616217
$ws.Range("M77").Value = 0
$ws.Range("N77").Value = 1327
$ws.Range("E78").Value = "-"
$ws.Range("E79").Value = 0
$ws.Range("F79").Value = -2529
$ws.Range("G79").Value = "-"
$ws.Range("E80").Value = -34816
$ws.Range("F80").Value = -41510
$ws.Range("G80").Value = -28419
$ws.Range("H80").Value = 28419
$ws.Range("I80").Value = 0
$ws.Range("E81").Value = 288488
$ws.Range("F81").Value = 726918
$ws.Range("G81").Value = 193933
$ws.Range("H81").Value = 310079
$ws.Range("I81").Value = 233686
$ws.Range("J81").Value = 224347
$ws.Range("K81").Value = 247321
$ws.Range("L81").Value = 236413
$ws.Range("M81").Value = 414348
$ws.Range("N81").Value = 428244
$ws.Range("E82").Value = 768149
$ws.Range("F82").Value = 1102129
$ws.Range("G82").Value = 668940
$ws.Range("H82").Value = 111749
$ws.Range("I82").Value = 248291
$ws.Range("J82").Value = -62583
$ws.Range("K82").Value = -57655
$ws.Range("L82").Value = 241880
$ws.Range("M82").Value = 414402
$ws.Range("N82").Value = 596812
$ws.Range("E83").Value = 4900
$ws.Range("F83").Value = 8231
$ws.Range("G83").Value = 2045
$ws.Range("H83").Value = 1779
$ws.Range("I83").Value = 82
$ws.Range("J83").Value = -2019
$ws.Range("K83").Value = 3057
$ws.Range("L83").Value = 3615
$ws.Range("M83").Value = 10711
$ws.Range("N83").Value = 39538
$ws.Range("E84").Value = 1026721
$ws.Range("F84").Value = 1793239
$ws.Range("G84").Value = 836499
$ws.Range("H84").Value = 452026
$ws.Range("I84").Value = 482059
$ws.Range("J84").Value = 159745
$ws.Range("K84").Value = 192723
$ws.Range("L84").Value = 481908
$ws.Range("M84").Value = 839461
$ws.Range("N84").Value = 1065921
